$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells keep their existing General/text display by
# forcing a Text number format before assignment (otherwise Excel would
# auto-convert numeric-looking strings like "283.03" or percentages like
# "1.71%" into numbers and silently drop significant trailing zeros).
$cellAddresses = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D28", "E28", "D29", "E29", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46")
foreach ($addr in $cellAddresses) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "283.03"
$ws.Range("E2").Value = "1.71%"
$ws.Range("D3").Value = "28.35"
$ws.Range("E3").Value = "4.06%"
$ws.Range("D4").Value = "5.026"
$ws.Range("E4").Value = "3.15%"
$ws.Range("D5").Value = "0.06501"
$ws.Range("E5").Value = "1.10%"
$ws.Range("D6").Value = "7.226"
$ws.Range("E6").Value = "3.18%"
$ws.Range("D7").Value = "1.381"
$ws.Range("E7").Value = "14.57%"
$ws.Range("D8").Value = "0.9183"
$ws.Range("E8").Value = "3.27%"
$ws.Range("D9").Value = "0.1542"
$ws.Range("E9").Value = "-0.59%"
$ws.Range("D10").Value = "0.06374"
$ws.Range("E10").Value = "24.91%"
$ws.Range("D11").Value = "0.07593"
$ws.Range("E11").Value = "1.26%"
$ws.Range("D12").Value = "0.02851"
$ws.Range("E12").Value = "-1.25%"
$ws.Range("D13").Value = "0.08976"
$ws.Range("E13").Value = "0.09%"
$ws.Range("D14").Value = "0.001589"
$ws.Range("E14").Value = "1.22%"
$ws.Range("D15").Value = "0.0006357"
$ws.Range("E15").Value = "-0.28%"
$ws.Range("D16").Value = "0.006124"
$ws.Range("E16").Value = "0.21%"
$ws.Range("E17").Value = "-0.97%"
$ws.Range("D18").Value = "3.360"
$ws.Range("E18").Value = "1.63%"
$ws.Range("D19").Value = "2.242"
$ws.Range("E19").Value = "1.07%"
$ws.Range("D20").Value = "0.3182"
$ws.Range("E20").Value = "-0.05%"
$ws.Range("E21").Value = "-2.21%"
$ws.Range("D22").Value = "3.968"
$ws.Range("E22").Value = "1.33%"
$ws.Range("E23").Value = "2.89%"
$ws.Range("D24").Value = "0.04442"
$ws.Range("E24").Value = "0.87%"
$ws.Range("D25").Value = "0.001184"
$ws.Range("E25").Value = "0.79%"
$ws.Range("D26").Value = "0.004458"
$ws.Range("E26").Value = "14.96%"
$ws.Range("D28").Value = "0.0001199"
$ws.Range("E28").Value = "1.61%"
$ws.Range("D29").Value = "0.0001618"
$ws.Range("E29").Value = "-1.58%"
$ws.Range("D40").Value = "0.04109"
$ws.Range("E40").Value = "-0.30%"
$ws.Range("D41").Value = "0.006662"
$ws.Range("E41").Value = "-2.06%"
$ws.Range("D42").Value = "0.1232"
$ws.Range("E42").Value = "4.93%"
$ws.Range("D43").Value = "0.002128"
$ws.Range("E43").Value = "10.85%"
$ws.Range("D44").Value = "0.01154"
$ws.Range("E44").Value = "-2.55%"
$ws.Range("D45").Value = "0.00005643"
$ws.Range("E45").Value = "5.91%"
$ws.Range("D46").Value = "1.954"
$ws.Range("E46").Value = "16.28%"
